$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2024-03-25 23:40:40"
$ws.Range("B2").Value = "https://wellfound.com/jobs/2963088-prompt-engineer-stem-freelancer"

$ws.Range("C2").Value = "Company Description:`nPhrasIQ is a leading AI company dedicated to advancing AI applications through the utilization of high-quality data. Our services harness enterprise data to customize robust base generative models, enabling the safe and efficient extraction of AI value. We provide tailored solutions for data collection, curation, and annotation, along with evaluation and optimization tools to ensure model efficacy. PhrasIQ empowers the development of cutting-edge LLMs and generative models globally through expertise in Reinforcement Learning from Human Feedback (RLHF), data generation, model evaluation, safety protocols, and alignment strategies.`nRole Description:`nAre you a subject matter expert in Computer Science, Mathematics, Advanced Physics, or Chemistry? PhrasIQ is actively seeking exceptional talent like you to play a pivotal role in training and enhancing Artificial Intelligence models, regardless of prior AI experience. Candidates holding Bachelor's, Master's, or PhDs in these fields are encouraged to apply.`nAs a Prompt Engineer, you will help the development and refinement of prompts for AI models, ensuring the quality and efficacy of the generated output. You'll collaborate closely with our dynamic team to design, implement, and evaluate prompt engineering strategies. Additionally, you'll stay abreast of the latest advancements in AI and STEM fields to ensure our models remain at the cutting edge of innovation.`nRoles and Responsibilities:`nEngage in flexible project-based roles, allowing you to tailor your workload to suit your schedule and interests.`nWork remotely, utilizing your own PC to contribute to our innovative projects.`nAssist in training advanced Large Language Models (LLMs), leveraging your expertise to refine and optimize model performance.`nAct as a domain expert in AI model fine-tuning, employing techniques such as Prompt Engineering and Reinforced Learning to enhance model capabilities.`nReview and generate technical content, providing valuable insights and expertise to enrich our AI models and deliver exceptional results.`nPay & Benefits:`nEnjoy project-based pay, with the potential to earn up to Rs. 2700 per hour based on the complexity and scope of your assignments.`nWith successful completion of our rigorous screening process, you could potentially earn up to Rs. 127K monthly, unlocking significant earning potential and rewarding your expertise and contributions.`nReady to shape the future of AI? Apply now and unleash your expertise with PhrasIQ!`nPhrasIQ is an equal opportunity employer. All qualified applicants will receive consideration for employment without regard to race, color, religion, gender, gender identity or expression, sexual orientation, national origin, genetics, disability, or age."

$ws.Range("D2").Value = "My expertise in Prompt Engineering and experience in AI model fine-tuning make me an ideal candidate for this role. I am eager to leverage my domain knowledge in Computer Science to enhance PhrasIQ's AI models.`n[Two of my projects relevant to the role are:`n1. [Synapse](https://github.com/UltigendLemate/synapse) - A productivity-enhancing collaborative web-app with features of Notion and Google Docs.`n2. [Quik Planr](https://github.com/UltigendLemate/Quik-Planr) - An application that turns users' two-liner ideas into comprehensive market research plans.]`nI am confident that my contributions can help PhrasIQ shape the future of AI."
